$wb = $excel.ActiveWorkbook

# Map of sheet name -> new value for cell B11
$updates = @{
    "Silver Rear_side"          = "5,289"
    "Silver Busbar front-side"  = "7,917"
    "Silver finger front-side"  = "7,967"
    "USD_CNY"                   = "7.2506"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cell = $ws.Range("B11")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$sheetName]
}
